# "Created experiment order generation script"
#
# The task-order workbook is regenerated: five sheets, reordered and
# renamed, each populated with a freshly generated set of stimulus-file
# task orders.
#
# Strategy: keep exactly one of the original sheets alive (so the
# workbook is never briefly empty), rename it first, then Copy() it
# four more times to build the other four sheets. Copy() (single-call
# "Copy(destination)" form, used consistently below) clones the
# row/column formatting exactly -- no new style entries get created --
# which keeps the bold/bordered header & index-column look intact.
# Afterwards we resize each sheet's row count (deleting or cloning
# rows via the same single-call Copy form) and overwrite the A/B
# column values with the new data.

$wb = $excel.ActiveWorkbook

# --- 1. Drop four of the five original sheets, keep one as the seed ---
$wb.Worksheets.Item("NB_TO-16512555599904544").Delete() | Out-Null
$wb.Worksheets.Item("RS_TO-16512555599914553").Delete() | Out-Null
$wb.Worksheets.Item("TOL_TO-16512555600554547").Delete() | Out-Null
$wb.Worksheets.Item("vSAT_TO-16512555601344533").Delete() | Out-Null

$seed = $wb.Worksheets.Item("GNG_TO-16512555566310518")
$seed.Name = "vSAT_TO-16515889935088773"

# --- 2. Rebuild the other four sheets by copying the seed (preserves styles) ---
$seed.Copy($null, $seed)
$sRS = $wb.Worksheets.Item(2)
$sRS.Name = "RS_TO-1651588993510725"

$sRS.Copy($null, $sRS)
$sGNG = $wb.Worksheets.Item(3)
$sGNG.Name = "GNG_TO-1651588993580242"

$sGNG.Copy($null, $sGNG)
$sNB = $wb.Worksheets.Item(4)
$sNB.Name = "NB_TO-16515889954322972"

$sNB.Copy($null, $sNB)
$sTOL = $wb.Worksheets.Item(5)
$sTOL.Name = "TOL_TO-16515889954791708"

$sVSAT = $wb.Worksheets.Item(1)

# --- 3. Resize each sheet to the row count its new data needs ---
# vSAT (sVSAT) and GNG (sGNG) already have 5 rows (header + 4) - no change.
# RS only needs 3 rows (header + 2) - drop rows 4:5.
$sRS.Rows("4:5").Delete()
# NB needs 10 rows (header + 9) - clone rows 6:9 from 2:5, then row 10 from row 5.
$sNB.Range("A2:B5").Copy($sNB.Range("A6:B9"))
$sNB.Range("A5:B5").Copy($sNB.Range("A10:B10"))
# TOL needs 7 rows (header + 6) - clone rows 6:7 from 4:5.
$sTOL.Range("A4:B5").Copy($sTOL.Range("A6:B7"))

# --- 4. Write the new cell values ---

# vSAT sheet
$sVSAT.Range("A2").Value = 0
$sVSAT.Range("B2").Value = "SAT_stims-16515889934346688.csv"
$sVSAT.Range("A3").Value = 1
$sVSAT.Range("B3").Value = "SAT_stims-16515889934476664.csv"
$sVSAT.Range("A4").Value = 2
$sVSAT.Range("B4").Value = "vSAT_stims-16515889934621298.csv"
$sVSAT.Range("A5").Value = 3
$sVSAT.Range("B5").Value = "vSAT_stims-16515889934771297.csv"

# RS sheet
$sRS.Range("A2").Value = 0
$sRS.Range("B2").Value = "eyes open"
$sRS.Range("A3").Value = 1
$sRS.Range("B3").Value = "eyes closed"

# GNG sheet
$sGNG.Range("A2").Value = 0
$sGNG.Range("B2").Value = "go_stims-1651588993516768.csv"
$sGNG.Range("A3").Value = 1
$sGNG.Range("B3").Value = "GNG_stims-16515889935405855.csv"
$sGNG.Range("A4").Value = 2
$sGNG.Range("B4").Value = "go_stims-1651588993548523.csv"
$sGNG.Range("A5").Value = 3
$sGNG.Range("B5").Value = "GNG_stims-1651588993580242.csv"

# NB sheet
$sNB.Range("A2").Value = 0
$sNB.Range("B2").Value = "OB-16515889941773844.csv"
$sNB.Range("A3").Value = 1
$sNB.Range("B3").Value = "ZB-match_6-16515889939621751.csv"
$sNB.Range("A4").Value = 2
$sNB.Range("B4").Value = "TB-16515889954166768.csv"
$sNB.Range("A5").Value = 3
$sNB.Range("B5").Value = "OB-16515889943531215.csv"
$sNB.Range("A6").Value = 4
$sNB.Range("B6").Value = "OB-1651588994613055.csv"
$sNB.Range("A7").Value = 5
$sNB.Range("B7").Value = "ZB-match_8-1651588993752906.csv"
$sNB.Range("A8").Value = 6
$sNB.Range("B8").Value = "TB-16515889947068317.csv"
$sNB.Range("A9").Value = 7
$sNB.Range("B9").Value = "ZB-match_3-16515889938033369.csv"
$sNB.Range("A10").Value = 8
$sNB.Range("B10").Value = "TB-16515889952581806.csv"

# TOL sheet
$sTOL.Range("A2").Value = 0
$sTOL.Range("B2").Value = "MM_stims-16515889954478898.csv"
$sTOL.Range("A3").Value = 1
$sTOL.Range("B3").Value = "ZM_stims-16515889954322972.csv"
$sTOL.Range("A4").Value = 2
$sTOL.Range("B4").Value = "MM_stims-16515889954635496.csv"
$sTOL.Range("A5").Value = 3
$sTOL.Range("B5").Value = "ZM_stims-16515889954478898.csv"
$sTOL.Range("A6").Value = 4
$sTOL.Range("B6").Value = "MM_stims-16515889954791708.csv"
$sTOL.Range("A7").Value = 5
$sTOL.Range("B7").Value = "ZM_stims-16515889954635496.csv"

# --- 5. Final sheet order: vSAT, RS, GNG, NB, TOL (already built in this order) ---
$sVSAT.Select()
